# Append " (Changed main)" after the existing sentence in the first
# paragraph, as three separate runs (" (", "Changed main", ")") that sit
# alongside the original, untouched run.
#
# A plain sequence of Range.InsertAfter calls would get silently
# coalesced back into the original run because the inserted text shares
# identical (empty) run formatting. To keep the runs distinct we turn on
# Track Changes for the duration of the three inserts - each tracked
# insertion is kept as its own revision/run - and then accept all
# revisions, which bakes the text in as plain runs while preserving the
# run boundaries that were created.

$d = $word.ActiveDocument

$wasTracking = $d.TrackRevisions
$d.TrackRevisions = $true

$r1 = $d.Paragraphs(1).Range
$r1.Collapse(0)      # wdCollapseEnd
$r1.MoveEnd(1, -1)   # step back before the paragraph mark
$r1.InsertAfter(" (")

$r2 = $d.Paragraphs(1).Range
$r2.Collapse(0)
$r2.MoveEnd(1, -1)
$r2.InsertAfter("Changed main")

$r3 = $d.Paragraphs(1).Range
$r3.Collapse(0)
$r3.MoveEnd(1, -1)
$r3.InsertAfter(")")

$d.TrackRevisions = $wasTracking
$d.AcceptAllRevisions()
